$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New predictor hypothesis label/complement columns for rows 5-9 (n_par = 2 block)
$ws.Range("E5").Value = "H1.V1>V2>0 "
$ws.Range("E6").Value = "H1.complement"
$ws.Range("E7").Value = "H2.V1>V2"
$ws.Range("E8").Value = "H2.complement"
$ws.Range("E9").Value = "Hu"

# Match formatting used by the existing E10:E16 "dimname" style block
$ws.Range("E10").Copy()
$ws.Range("E9").PasteSpecial(-4122)

# Extra empty formatted cell introduced alongside the new columns
$ws.Range("C5").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$null = $ws.Range("H11").ClearContents()

$excel.CutCopyMode = 0

# Reset view: scroll back to top-left and move the selection
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("E6").Select()
